# Update NATMI TPM-derived metrics (ligand/receptor specificity columns G:T)
# for each data row, per the refreshed TPM computation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 6.314527000000001
$ws.Cells.Item(2, 8).Value = 18.943581
$ws.Cells.Item(2, 9).Value = 0.2616724966426195
$ws.Cells.Item(2, 10).Value = 0.2616724966426195
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 0.239548
$ws.Cells.Item(2, 14).Value = 0.7186440000000001
$ws.Cells.Item(2, 15).Value = 0.2072176292017679
$ws.Cells.Item(2, 16).Value = 0.2072176292017679
$ws.Cells.Item(2, 17).Value = 1.512632313796
$ws.Cells.Item(2, 18).Value = 13.613690824164
$ws.Cells.Item(2, 19).Value = 0.0542231543815912
$ws.Cells.Item(2, 20).Value = 0.0542231543815912

# Row 3
$ws.Cells.Item(3, 7).Value = 6.314527000000001
$ws.Cells.Item(3, 8).Value = 18.943581
$ws.Cells.Item(3, 9).Value = 0.2616724966426195
$ws.Cells.Item(3, 10).Value = 0.2616724966426195
$ws.Cells.Item(3, 15).Value = 0.1375767575223525
$ws.Cells.Item(3, 16).Value = 0.1375767575223525
$ws.Cells.Item(3, 17).Value = 1.004272898291667
$ws.Cells.Item(3, 18).Value = 9.038456084625
$ws.Cells.Item(3, 19).Value = 0.03600005362087027
$ws.Cells.Item(3, 20).Value = 0.03600005362087028

# Row 4
$ws.Cells.Item(4, 7).Value = 6.314527000000001
$ws.Cells.Item(4, 8).Value = 18.943581
$ws.Cells.Item(4, 9).Value = 0.2616724966426195
$ws.Cells.Item(4, 10).Value = 0.2616724966426195
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 0.5662826666666666
$ws.Cells.Item(4, 14).Value = 1.698848
$ws.Cells.Item(4, 15).Value = 0.4898548585031879
$ws.Cells.Item(4, 16).Value = 0.4898548585031879
$ws.Cells.Item(4, 17).Value = 3.575807188298667
$ws.Cells.Item(4, 18).Value = 32.182264694688
$ws.Cells.Item(4, 19).Value = 0.1281815438170463
$ws.Cells.Item(4, 20).Value = 0.1281815438170463

# Row 5
$ws.Cells.Item(5, 7).Value = 6.314527000000001
$ws.Cells.Item(5, 8).Value = 18.943581
$ws.Cells.Item(5, 9).Value = 0.2616724966426195
$ws.Cells.Item(5, 10).Value = 0.2616724966426195
$ws.Cells.Item(5, 13).Value = 0.191149
$ws.Cells.Item(5, 14).Value = 0.573447
$ws.Cells.Item(5, 15).Value = 0.1653507547726916
$ws.Cells.Item(5, 16).Value = 0.1653507547726916
$ws.Cells.Item(5, 17).Value = 1.207015521523
$ws.Cells.Item(5, 18).Value = 10.863139693707
$ws.Cells.Item(5, 19).Value = 0.04326774482311176
$ws.Cells.Item(5, 20).Value = 0.04326774482311176

# Row 6
$ws.Cells.Item(6, 9).Value = 0.1461016137776048
$ws.Cells.Item(6, 10).Value = 0.1461016137776048
$ws.Cells.Item(6, 11).Value = 2
$ws.Cells.Item(6, 12).Value = 0.6666666666666666
$ws.Cells.Item(6, 13).Value = 0.239548
$ws.Cells.Item(6, 14).Value = 0.7186440000000001
$ws.Cells.Item(6, 15).Value = 0.2072176292017679
$ws.Cells.Item(6, 16).Value = 0.2072176292017679
$ws.Cells.Item(6, 17).Value = 0.8445596114733334
$ws.Cells.Item(6, 18).Value = 7.60103650326
$ws.Cells.Item(6, 19).Value = 0.03027483002954761
$ws.Cells.Item(6, 20).Value = 0.03027483002954761

# Row 7
$ws.Cells.Item(7, 9).Value = 0.1461016137776048
$ws.Cells.Item(7, 10).Value = 0.1461016137776048
$ws.Cells.Item(7, 15).Value = 0.1375767575223525
$ws.Cells.Item(7, 16).Value = 0.1375767575223525
$ws.Cells.Item(7, 19).Value = 0.02010018629230593
$ws.Cells.Item(7, 20).Value = 0.02010018629230594

# Row 8
$ws.Cells.Item(8, 9).Value = 0.1461016137776048
$ws.Cells.Item(8, 10).Value = 0.1461016137776048
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 0.5662826666666666
$ws.Cells.Item(8, 14).Value = 1.698848
$ws.Cells.Item(8, 15).Value = 0.4898548585031879
$ws.Cells.Item(8, 16).Value = 0.4898548585031879
$ws.Cells.Item(8, 17).Value = 1.996507877102222
$ws.Cells.Item(8, 18).Value = 17.96857089392
$ws.Cells.Item(8, 19).Value = 0.07156858534411599
$ws.Cells.Item(8, 20).Value = 0.071568585344116

# Row 9
$ws.Cells.Item(9, 9).Value = 0.1461016137776048
$ws.Cells.Item(9, 10).Value = 0.1461016137776048
$ws.Cells.Item(9, 13).Value = 0.191149
$ws.Cells.Item(9, 14).Value = 0.573447
$ws.Cells.Item(9, 15).Value = 0.1653507547726916
$ws.Cells.Item(9, 16).Value = 0.1653507547726916
$ws.Cells.Item(9, 17).Value = 0.6739222417783334
$ws.Cells.Item(9, 18).Value = 6.065300176005
$ws.Cells.Item(9, 19).Value = 0.02415801211163523
$ws.Cells.Item(9, 20).Value = 0.02415801211163523

# Row 10
$ws.Cells.Item(10, 7).Value = 0.510814
$ws.Cells.Item(10, 8).Value = 1.532442
$ws.Cells.Item(10, 9).Value = 0.02116801063642661
$ws.Cells.Item(10, 10).Value = 0.02116801063642662
$ws.Cells.Item(10, 11).Value = 2
$ws.Cells.Item(10, 12).Value = 0.6666666666666666
$ws.Cells.Item(10, 13).Value = 0.239548
$ws.Cells.Item(10, 14).Value = 0.7186440000000001
$ws.Cells.Item(10, 15).Value = 0.2072176292017679
$ws.Cells.Item(10, 16).Value = 0.2072176292017679
$ws.Cells.Item(10, 17).Value = 0.122364472072
$ws.Cells.Item(10, 18).Value = 1.101280248648
$ws.Cells.Item(10, 19).Value = 0.004386384978998129
$ws.Cells.Item(10, 20).Value = 0.00438638497899813

# Row 11
$ws.Cells.Item(11, 7).Value = 0.510814
$ws.Cells.Item(11, 8).Value = 1.532442
$ws.Cells.Item(11, 9).Value = 0.02116801063642661
$ws.Cells.Item(11, 10).Value = 0.02116801063642662
$ws.Cells.Item(11, 15).Value = 0.1375767575223525
$ws.Cells.Item(11, 16).Value = 0.1375767575223525
$ws.Cells.Item(11, 17).Value = 0.08124070991666667
$ws.Cells.Item(11, 18).Value = 0.7311663892500001
$ws.Cells.Item(11, 19).Value = 0.002912226266558243
$ws.Cells.Item(11, 20).Value = 0.002912226266558244

# Row 12
$ws.Cells.Item(12, 7).Value = 0.510814
$ws.Cells.Item(12, 8).Value = 1.532442
$ws.Cells.Item(12, 9).Value = 0.02116801063642661
$ws.Cells.Item(12, 10).Value = 0.02116801063642662
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 0.5662826666666666
$ws.Cells.Item(12, 14).Value = 1.698848
$ws.Cells.Item(12, 15).Value = 0.4898548585031879
$ws.Cells.Item(12, 16).Value = 0.4898548585031879
$ws.Cells.Item(12, 17).Value = 0.2892651140906666
$ws.Cells.Item(12, 18).Value = 2.603386026816
$ws.Cells.Item(12, 19).Value = 0.01036925285510074
$ws.Cells.Item(12, 20).Value = 0.01036925285510074

# Row 13
$ws.Cells.Item(13, 7).Value = 0.510814
$ws.Cells.Item(13, 8).Value = 1.532442
$ws.Cells.Item(13, 9).Value = 0.02116801063642661
$ws.Cells.Item(13, 10).Value = 0.02116801063642662
$ws.Cells.Item(13, 13).Value = 0.191149
$ws.Cells.Item(13, 14).Value = 0.573447
$ws.Cells.Item(13, 15).Value = 0.1653507547726916
$ws.Cells.Item(13, 16).Value = 0.1653507547726916
$ws.Cells.Item(13, 17).Value = 0.097641585286
$ws.Cells.Item(13, 18).Value = 0.8787742675740001
$ws.Cells.Item(13, 19).Value = 0.003500146535769505
$ws.Cells.Item(13, 20).Value = 0.003500146535769506

# Row 14
$ws.Cells.Item(14, 7).Value = 13.78043333333333
$ws.Cells.Item(14, 8).Value = 41.3413
$ws.Cells.Item(14, 9).Value = 0.571057878943349
$ws.Cells.Item(14, 10).Value = 0.5710578789433491
$ws.Cells.Item(14, 11).Value = 2
$ws.Cells.Item(14, 12).Value = 0.6666666666666666
$ws.Cells.Item(14, 13).Value = 0.239548
$ws.Cells.Item(14, 14).Value = 0.7186440000000001
$ws.Cells.Item(14, 15).Value = 0.2072176292017679
$ws.Cells.Item(14, 16).Value = 0.2072176292017679
$ws.Cells.Item(14, 17).Value = 3.301075244133334
$ws.Cells.Item(14, 18).Value = 29.7096771972
$ws.Cells.Item(14, 19).Value = 0.118333259811631
$ws.Cells.Item(14, 20).Value = 0.118333259811631

# Row 15
$ws.Cells.Item(15, 7).Value = 13.78043333333333
$ws.Cells.Item(15, 8).Value = 41.3413
$ws.Cells.Item(15, 9).Value = 0.571057878943349
$ws.Cells.Item(15, 10).Value = 0.5710578789433491
$ws.Cells.Item(15, 15).Value = 0.1375767575223525
$ws.Cells.Item(15, 16).Value = 0.1375767575223525
$ws.Cells.Item(15, 17).Value = 2.191663084722222
$ws.Cells.Item(15, 18).Value = 19.7249677625
$ws.Cells.Item(15, 19).Value = 0.07856429134261807
$ws.Cells.Item(15, 20).Value = 0.07856429134261809

# Row 16
$ws.Cells.Item(16, 7).Value = 13.78043333333333
$ws.Cells.Item(16, 8).Value = 41.3413
$ws.Cells.Item(16, 9).Value = 0.571057878943349
$ws.Cells.Item(16, 10).Value = 0.5710578789433491
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 0.5662826666666666
$ws.Cells.Item(16, 14).Value = 1.698848
$ws.Cells.Item(16, 15).Value = 0.4898548585031879
$ws.Cells.Item(16, 16).Value = 0.4898548585031879
$ws.Cells.Item(16, 17).Value = 7.803620535822222
$ws.Cells.Item(16, 18).Value = 70.2325848224
$ws.Cells.Item(16, 19).Value = 0.2797354764869248
$ws.Cells.Item(16, 20).Value = 0.2797354764869249

# Row 17
$ws.Cells.Item(17, 7).Value = 13.78043333333333
$ws.Cells.Item(17, 8).Value = 41.3413
$ws.Cells.Item(17, 9).Value = 0.571057878943349
$ws.Cells.Item(17, 10).Value = 0.5710578789433491
$ws.Cells.Item(17, 13).Value = 0.191149
$ws.Cells.Item(17, 14).Value = 0.573447
$ws.Cells.Item(17, 15).Value = 0.1653507547726916
$ws.Cells.Item(17, 16).Value = 0.1653507547726916
$ws.Cells.Item(17, 17).Value = 2.634116051233334
$ws.Cells.Item(17, 18).Value = 23.7070444611
$ws.Cells.Item(17, 19).Value = 0.09442485130217512
$ws.Cells.Item(17, 20).Value = 0.09442485130217514
